$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testDataAPI")

# Replace the existing customer id in A13 with a new one, and add another new id below it
$ws.Range("A13").Value = "cus_JQcOWxkfkURlhm"
$ws.Range("A14").Value = "cus_JQcNaXiIYdviRe"

# Update the email in B9 (the "at12@zmail.com" contact gets the "3" appended)
$ws.Range("B9").Value = "at123@zmail.com"

# Match the saved selection state from the authored workbook
$ws.Range("B9").Select()
